$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet,
#    re-using the formatting of the existing quarter sheets.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy header / column-A styling (bold, centred, bordered) from the
# "2021-Q4" sheet, which has the same 8-column fund-holding layout.
$q4.Range("A1:H7").Copy()
$q1.Range("A1").PasteSpecial(-4122)

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

$q1Rows = @(
    @(0, "006022", "富国大盘价值量化精选混合", "3.47", "86.62", "1.35", "0.0468", 9),
    @(1, "519616", "银河君信灵活配置混合A",   "4.54", "24.42", "0.85", "0.0386", 5),
    @(2, "519618", "银河君信灵活配置混合I",   "4.54", "24.42", "0.85", "0.0386", 5),
    @(3, "519656", "银河灵活配置混合 - A",    "0.72", "59.27", "4.04", "0.0291", 3),
    @(4, "519657", "银河灵活配置混合 - C",    "0.33", "59.27", "4.04", "0.0133", 3),
    @(5, "519617", "银河君信灵活配置混合C",   "0.64", "24.42", "0.85", "0.0054", 5)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r,1).Value = $row[0]
    $q1.Cells.Item($r,2).NumberFormat = "@"
    $q1.Cells.Item($r,2).Value = $row[1]
    $q1.Cells.Item($r,3).Value = $row[2]
    $q1.Cells.Item($r,4).NumberFormat = "@"
    $q1.Cells.Item($r,4).Value = $row[3]
    $q1.Cells.Item($r,5).NumberFormat = "@"
    $q1.Cells.Item($r,5).Value = $row[4]
    $q1.Cells.Item($r,6).NumberFormat = "@"
    $q1.Cells.Item($r,6).Value = $row[5]
    $q1.Cells.Item($r,7).NumberFormat = "@"
    $q1.Cells.Item($r,7).Value = $row[6]
    $q1.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q1, pushing the existing 2021-Q4 / 2021-Q3 rows down.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# New row takes the same formatting as the (now-shifted) row below it.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 0.17

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2

# ------------------------------------------------------------------
# 3. Restore the originally-active sheet so this edit doesn't
#    incidentally change which tab is shown on re-open.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$q3.Activate()

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
